$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 6.15379541431027

# Row 3 values
$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 0.1496068669990043
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 4.260211312413533
